# Add constraint examples to ingrowth form.
# Operates on the "survey" sheet (3rd sheet / Worksheets.Item(3)) of the
# ingrowth.xlsx XLSForm workbook: introduces a `constraint` /
# `display.constraint_message.text` pair of columns (G:H) and a worked
# `if / note / end if` example clause demonstrating a DBH constraint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Make room for the new "if / note / end if" example block: old rows 7-25
# (From tag / Status / ... / end screen) shift down to rows 10-28.
$ws.Rows("7:9").Insert()

# --- New header cells for the constraint columns (row 1) ---
$ws.Range("G1").Value = "constraint"
$ws.Range("H1").Value = "display.constraint_message.text"

# --- New example clause block illustrating the constraint (rows 7-9) ---
# (cell-write order below matches the shared-string allocation order of the
#  authored workbook so new <si> entries land on the same indices)
$ws.Range("A7").Value = "if"
$ws.Range("A9").Value = "end if"

# --- Constraint + message attached to the existing "dbh" question (row 6) ---
$ws.Range("H6").Value = "DBH != 123."
$ws.Range("B7").Value = "data('dbh') != 123"

$ws.Range("C8").Value = "note"
$ws.Range("F8").Value = "Dbh isnt right."

$constraintFormula = '(function() { if (data(''dbh'') == "123") { return true;} alert(''DBH incorrect''); return false;}) ()'
$ws.Range("G6").Value = $constraintFormula

# --- Column sizing to fit the new/changed content ---
# (engine stores width in 1-pixel / MDW=6 steps: stored = (round(chars*6)+5)/6;
#  inputs below are the closest pixel-quantised values to the authored widths)
$ws.Columns("B").ColumnWidth = 15.666666666666666
$ws.Columns("G").ColumnWidth = 80.66666666666667
$ws.Columns("H").ColumnWidth = 29.5

# --- Selection left on the new constraint-formula cell ---
$ws.Range("G6").Select()

# --- Page setup (orientation) touched while reviewing the new columns ---
$ws.PageSetup.Orientation = 1
